# Update "想去人数" (column F) figures on the "展览" and "全部类型" sheets
# to reflect newly generated output counts.

$wb = $excel.ActiveWorkbook

# Map of row -> new value for sheet "展览" (sheet index 1)
$exhibitionUpdates = @{
    2  = 15105
    3  = 19321
    14 = 191
    15 = 237
    17 = 1495
    21 = 242
    22 = 8095
    27 = 1261
    30 = 6103
    31 = 125
    33 = 179
    35 = 296
    36 = 5511
    37 = 1009
    38 = 23
    39 = 31
    40 = 56
}

# Map of row -> new value for sheet "全部类型" (sheet index 4)
$allTypesUpdates = @{
    2  = 15105
    3  = 19321
    14 = 191
    15 = 237
    17 = 1495
    22 = 242
    23 = 8095
    28 = 1261
    33 = 6103
    34 = 125
    36 = 179
    38 = 296
    39 = 5511
    40 = 1009
    41 = 23
    42 = 31
    43 = 56
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Range("F$row").Value = $allTypesUpdates[$row]
}
